$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '37.095.71'
$ws.Range('E2').Value = '  -0.14%  '
$ws.Range('D3').Value = '2.056.44'
$ws.Range('E3').Value = '  -1.29%  '
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').Value = "'249.58"
$ws.Range('E5').Value = '  -0.88%  '
$ws.Range('D6').Value = "'0.673"
$ws.Range('E6').Value = '  +2.11%  '
$ws.Range('D7').Value = "'59.06"
$ws.Range('E7').Value = '  +15.61%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').Value = "'60.82"
$ws.Range('E9').Value = '  -0.12%  '
$ws.Range('D10').Value = "'0.378"
$ws.Range('E10').Value = '  +1.13%  '
$ws.Range('D11').Value = "'0.0803"
$ws.Range('E11').Value = '  +7.54%  '
$ws.Range('E12').Value = '  +0.85%  '
$ws.Range('D13').Value = "'15.14"
$ws.Range('E13').Value = '  +1.57%  '
$ws.Range('E14').Value = '  -1.24%  '
$ws.Range('D15').Value = "'0.816"
$ws.Range('E15').Value = '  -2.34%  '
$ws.Range('D16').Value = "'5.31"
$ws.Range('E16').Value = '  +3.56%  '
$ws.Range('D17').Value = '2.064.74'
$ws.Range('E17').Value = '  -1.02%  '
$ws.Range('D18').Value = '37.117.52'
$ws.Range('E18').Value = '  +0.19%  '
$ws.Range('D19').Value = "'74.65"
$ws.Range('E19').Value = '  +2.91%  '
$ws.Range('D20').Value = '0.0₃0919'
$ws.Range('E20').Value = '  +11.09%  '
$ws.Range('D21').Value = "'14.49"
$ws.Range('E21').Value = '  +8.02%  '
$ws.Range('D22').Value = "'5.34"
$ws.Range('E22').Value = '  +2.08%  '
$ws.Range('D23').Value = "'238.44"
$ws.Range('E23').Value = '  -0.87%  '
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('D25').Value = "'2.43"
$ws.Range('E25').Value = '  -2.59%  '
$ws.Range('D26').Value = "'171.83"
$ws.Range('E26').Value = '  +0.81%  '
$ws.Range('E27').Value = '  -1.61%  '
$ws.Range('D28').Value = "'20.18"
$ws.Range('E28').Value = '  -4.16%  '
$ws.Range('E29').Value = '  -0.48%  '
$ws.Range('E30').Value = '  +1.82%  '
$ws.Range('E31').Value = '  -0.73%  '
$ws.Range('D32').Value = "'4.61"
$ws.Range('E32').Value = '  +2.36%  '
$ws.Range('D33').Value = "'0.0633"
$ws.Range('E33').Value = '  +3.78%  '
$ws.Range('E34').Value = '  +7.41%  '
$ws.Range('D35').Value = "'0.0881"
$ws.Range('E35').Value = '  -5.63%  '
$ws.Range('E36').Value = '  -0.06%  '
$ws.Range('D37').Value = "'2.27"
$ws.Range('E37').Value = '  -1.60%  '
$ws.Range('E38').Value = '  -2.51%  '
$ws.Range('D39').Value = "'0.109"
$ws.Range('E39').Value = '  +25.64%  '
$ws.Range('D40').Value = "'1.35"
$ws.Range('E40').Value = '  +1.56%  '
$ws.Range('D41').Value = "'18.30"
$ws.Range('E41').Value = '  +3.30%  '
$ws.Range('E42').Value = '  +0.41%  '
$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').Value = "'1.14"
$ws.Range('E43').Value = '  -0.77%  '
$ws.Range('B44').Value = 'FTXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D44').Value = "'4.43"
$ws.Range('E44').Value = '  +37.81%  '
$ws.Range('D45').Value = "'96.91"
$ws.Range('E45').Value = '  -1.07%  '
$ws.Range('E46').Value = '  -0.63%  '
$ws.Range('D47').Value = "'4.47"
$ws.Range('E47').Value = '  +13.17%  '
$ws.Range('E48').Value = '  +9.15%  '
$ws.Range('D49').Value = '1.302.93'
$ws.Range('E49').Value = '  -1.30%  '
$ws.Range('D50').Value = "'2.91"
$ws.Range('E50').Value = '  -1.40%  '
$ws.Range('D51').Value = "'6.89"
$ws.Range('E51').Value = '  -0.90%  '
